# This script expands every shift that was shared by multiple people
# (e.g. "Jack Green, Daniel Senn") into one row per assigned person, and
# refreshes the "Full Schedule" summary sheet to match.
#
# It rewrites each day sheet's data rows (everything below the header row)
# and then rebuilds the Full Schedule sheet as the concatenation of all the
# day sheets, in the same day order used by that sheet originally.

$wb = $excel.ActiveWorkbook

# Target data (Start, End, Assigned) for each day sheet, after splitting
# combined "Name1, Name2" assignments into separate rows.
$sunday = @(
    @("9:00 AM", "12:00 PM", "Daniel Senn"),
    @("9:00 AM", "12:00 PM", "Mario Sell"),
    @("12:00 PM", "2:00 PM", "Brooke Kazmierczak"),
    @("12:00 PM", "2:00 PM", "Patrick Testing"),
    @("2:00 PM", "5:00 PM", "Zion Williams"),
    @("2:00 PM", "5:00 PM", "Brooke Kazmierczak")
)

$monday = @(
    @("9:00 AM", "2:00 PM", "Tatiana Mata Diaz"),
    @("12:00 PM", "5:00 PM", "Billy Happy")
)

$tuesday = @(
    @("9:00 AM", "2:00 PM", "Rick Kemper"),
    @("2:00 PM", "5:00 PM", "Jash Hitesh Parekh")
)

$wednesday = @(
    @("9:00 AM", "1:00 PM", "Alan Haim"),
    @("9:00 AM", "1:00 PM", "Thomas Mack"),
    @("1:00 PM", "5:00 PM", "Daniel Finn"),
    @("1:00 PM", "5:00 PM", "Nikko Sandgren")
)

$thursday = @(
    @("2:00 PM", "4:00 PM", "Jash Hitesh Parekh"),
    @("9:00 AM", "12:00 PM", "Mario Sell"),
    @("9:00 AM", "12:00 PM", "Matthew Chase"),
    @("12:00 PM", "2:00 PM", "Gissel O Rosa"),
    @("12:00 PM", "2:00 PM", "Chris Test")
)

$friday = @(
    @("9:00 AM", "12:00 PM", "Todd chop"),
    @("9:00 AM", "12:00 PM", "Sebastian Hurd"),
    @("12:00 PM", "2:00 PM", "Brooke Kazmierczak"),
    @("12:00 PM", "2:00 PM", "Jullian Kemp"),
    @("2:00 PM", "5:00 PM", "Olivia Schindler"),
    @("2:00 PM", "5:00 PM", "Regenae Walkters")
)

$saturday = @(
    @("9:00 AM", "2:00 PM", "Jack Green"),
    @("9:00 AM", "2:00 PM", "Krish Chawla"),
    @("2:00 PM", "4:00 PM", "Patrick Testing"),
    @("2:00 PM", "4:00 PM", "Daniel Senn"),
    @("4:00 PM", "5:00 PM", "Zion Williams"),
    @("4:00 PM", "5:00 PM", "Greg Aivaliotis")
)

function Update-DaySheet($sheetName, $rows) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove all existing data (everything below the header row) so the
    # sheet's used range shrinks back down before we repopulate it.
    $used = $ws.UsedRange
    $lastRow = $used.Row + $used.Rows.Count - 1
    if ($lastRow -ge 2) {
        $ws.Range("A2:C" + $lastRow).ClearContents()
    }

    $r = 2
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $r = $r + 1
    }
}

Update-DaySheet "Sunday" $sunday
Update-DaySheet "Monday" $monday
Update-DaySheet "Tuesday" $tuesday
Update-DaySheet "Wednesday" $wednesday
Update-DaySheet "Thursday" $thursday
Update-DaySheet "Friday" $friday
Update-DaySheet "Saturday" $saturday

# Rebuild the "Full Schedule" sheet as the union of all day sheets, using
# the same day ordering as the original sheet (Monday, Tuesday, Thursday,
# Wednesday, Friday, Saturday, Sunday).
$full = $wb.Worksheets.Item("Full Schedule")
$usedFull = $full.UsedRange
$lastFullRow = $usedFull.Row + $usedFull.Rows.Count - 1
if ($lastFullRow -ge 2) {
    $full.Range("A2:D" + $lastFullRow).ClearContents()
}

$dayOrder = @(
    @("Monday", $monday),
    @("Tuesday", $tuesday),
    @("Thursday", $thursday),
    @("Wednesday", $wednesday),
    @("Friday", $friday),
    @("Saturday", $saturday),
    @("Sunday", $sunday)
)

$fr = 2
foreach ($entry in $dayOrder) {
    $dayName = $entry[0]
    $dayRows = $entry[1]
    foreach ($row in $dayRows) {
        $full.Cells.Item($fr, 1).Value = $dayName
        $full.Cells.Item($fr, 2).Value = $row[0]
        $full.Cells.Item($fr, 3).Value = $row[1]
        $full.Cells.Item($fr, 4).Value = $row[2]
        $fr = $fr + 1
    }
}

Write-Host "Schedule updated."
